$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 214, shifting rows 214:320 down to 215:321
$ws.Rows("214:214").Insert()

# Populate the newly inserted row 214 with the new weekly price entry.
$ws.Range("A214").Value = 5
$ws.Range("B214").Value = "Macroferia Regional de Talca"
$ws.Range("C214").Value = "Maule"
$ws.Range("D214").Value = 44452
$ws.Range("E214").Value = 7
$ws.Range("F214").Value = 100112020
$ws.Range("G214").Value = "Tomate"
$ws.Range("H214").Value = "Larga vida"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 1500
$ws.Range("K214").Value = 10000
$ws.Range("L214").Value = 10000
$ws.Range("M214").Value = 10000
$ws.Range("N214").Value = "`$/caja 10 kilos"
$ws.Range("O214").Value = "Región de Arica y Parinacota"
$ws.Range("P214").Value = 1000
$ws.Range("Q214").Value = 10
$ws.Range("R214").Value = "Hortaliza"
